$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (unchanged content, but keep for clarity)
$ws.Range("A1").Value = "data_bases"
$ws.Range("B1").Value = "begins"
$ws.Range("C1").Value = "ends"

# Update data rows with new GSE identifiers and new "ends" value
$ws.Range("A2").Value = "GSE43414"
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 120

$ws.Range("A3").Value = "GSE88890"
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 120

$ws.Range("A4").Value = "GSE41826"
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 120

# Remove the now-extra 5th row entirely
$ws.Rows.Item(5).Delete()

# Update selection to reflect the new active cell
$ws.Range("C4").Select()

# Match the author's print setup (portrait orientation)
$ws.PageSetup.Orientation = 1
